$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3513.889
$ws.Range("I19").Value = 3315.6
$ws.Range("J19").Value = 3761.75
$ws.Range("K19").Value = 3315.6
$ws.Range("L19").Value = 3761.75
$ws.Range("M19").Value = -3140.6
$ws.Range("N19").Value = -4111.75

$ws.Range("H64").Value = 9999.272000000001
$ws.Range("I64").Value = 5714
$ws.Range("J64").Value = 17498.5
$ws.Range("K64").Value = 5714
$ws.Range("L64").Value = 17498.5
$ws.Range("M64").Value = -5466
$ws.Range("N64").Value = -17994.5

$ws.Range("H67").Value = 9999.272000000001
$ws.Range("I67").Value = 5714
$ws.Range("J67").Value = 17498.5
$ws.Range("K67").Value = 5714
$ws.Range("L67").Value = 17498.5
$ws.Range("M67").Value = -4856
$ws.Range("N67").Value = -19214.5

$ws.Range("H98").Value = 1019.8889
$ws.Range("I98").Value = 1019.8889
$ws.Range("K98").Value = 1019.8889
$ws.Range("M98").Value = 478.1111

$ws.Range("H122").Value = 1019.8889
$ws.Range("I122").Value = 1019.8889
$ws.Range("K122").Value = 3059.6667
$ws.Range("M122").Value = -609.6667000000002

$ws.Range("H131").Value = 1633.3334
$ws.Range("I131").Value = 1633.3334
$ws.Range("K131").Value = 4900.0002
$ws.Range("M131").Value = 139.9997999999996

$ws.Range("H132").Value = 2701.6086
$ws.Range("I132").Value = 2698.7334
$ws.Range("K132").Value = 8096.2002
$ws.Range("M132").Value = -5566.2002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3094
$ws.Range("I45").Value = 1455.75
$ws.Range("K45").Value = 1455.75
$ws.Range("M45").Value = -1078.75

$ws.Range("H63").Value = 4506.8
$ws.Range("I63").Value = 2314.6
$ws.Range("J63").Value = 6699
$ws.Range("K63").Value = 2314.6
$ws.Range("L63").Value = 6699
$ws.Range("M63").Value = -1628.6
$ws.Range("N63").Value = -8071

$ws.Range("H66").Value = 4506.8
$ws.Range("I66").Value = 2314.6
$ws.Range("J66").Value = 6699
$ws.Range("K66").Value = 11573
$ws.Range("L66").Value = 33495
$ws.Range("M66").Value = -8141
$ws.Range("N66").Value = -40359

$ws.Range("H74").Value = 990.625
$ws.Range("I74").Value = 1007.9545
$ws.Range("K74").Value = 1007.9545
$ws.Range("M74").Value = -133.9545000000001

$ws.Range("H77").Value = 990.625
$ws.Range("I77").Value = 1007.9545
$ws.Range("K77").Value = 5039.7725
$ws.Range("M77").Value = -671.7725

$ws.Range("H97").Value = 1027.1666
$ws.Range("I97").Value = 632.6
$ws.Range("J97").Value = 3000
$ws.Range("K97").Value = 632.6
$ws.Range("L97").Value = 3000
$ws.Range("M97").Value = -136.6
$ws.Range("N97").Value = -3992

$ws.Range("H102").Value = 1000
$ws.Range("I102").Value = 1000
$ws.Range("K102").Value = 1000
$ws.Range("M102").Value = 622

$ws.Range("H121").Value = 47751.332
$ws.Range("J121").Value = 47751.332
$ws.Range("L121").Value = 47751.332
$ws.Range("N121").Value = -51245.332

$ws.Range("H122").Value = 2822.7368
$ws.Range("I122").Value = 2588.8462
$ws.Range("K122").Value = 7766.5386
$ws.Range("M122").Value = -5316.5386

$ws.Range("H132").Value = 1124.6364
$ws.Range("J132").Value = 900
$ws.Range("L132").Value = 2700
$ws.Range("N132").Value = -7760

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 248
$ws.Range("I11").Value = 248
$ws.Range("K11").Value = 248
$ws.Range("M11").Value = -108

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 5000
$ws.Range("L4").Value = 5000
$ws.Range("N4").Value = -5224

$ws.Range("H22").Value = 725
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()

$ws.Range("H94").Value = 3726.7778
$ws.Range("J94").Value = 3842.2
$ws.Range("L94").Value = 3842.2
$ws.Range("N94").Value = -4744.2

$ws.Range("H109").Value = 16600
$ws.Range("J109").Value = 17000
$ws.Range("L109").Value = 17000
$ws.Range("N109").Value = -19080

$ws.Range("H134").Value = 2256.2727
$ws.Range("I134").Value = 1958.1177
$ws.Range("J134").Value = 3270
$ws.Range("K134").Value = 5874.3531
$ws.Range("L134").Value = 9810
$ws.Range("M134").Value = -3339.3531
$ws.Range("N134").Value = -14880

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 368
$ws.Range("I8").Value = 368
$ws.Range("K8").Value = 1104
$ws.Range("M8").Value = -965

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws.Range("H105").Value = 12167.5
$ws.Range("J105").Value = 12167.5
$ws.Range("L105").Value = 12167.5
$ws.Range("N105").Value = -19155.5

$ws.Range("H113").Value = 2034.7142
$ws.Range("I113").Value = 1848.8
$ws.Range("J113").Value = 2499.5
$ws.Range("K113").Value = 1848.8
$ws.Range("L113").Value = 2499.5
$ws.Range("M113").Value = 321.2
$ws.Range("N113").Value = -6839.5

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -54920

$ws.Range("H132").Value = 8142.4287
$ws.Range("I132").Value = 8142.4287
$ws.Range("K132").Value = 24427.2861
$ws.Range("M132").Value = -21897.2861

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 1214.091
$ws.Range("I9").Value = 380.5
$ws.Range("J9").Value = 3437
$ws.Range("K9").Value = 380.5
$ws.Range("L9").Value = 3437
$ws.Range("M9").Value = -156.5
$ws.Range("N9").Value = -3885

$ws.Range("H40").Value = 4136
$ws.Range("I40").Value = 3806
$ws.Range("J40").Value = 5126
$ws.Range("K40").Value = 3806
$ws.Range("L40").Value = 5126
$ws.Range("M40").Value = -3670
$ws.Range("N40").Value = -5398

$ws.Range("H100").Value = 4743
$ws.Range("I100").Value = 3691.6
$ws.Range("K100").Value = 3691.6
$ws.Range("M100").Value = -3150.6

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").ClearContents()

$ws.Range("H132").Value = 2646.0386
$ws.Range("I132").Value = 2580.8572
$ws.Range("K132").Value = 7742.571599999999
$ws.Range("M132").Value = -5212.571599999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3499.6667
$ws.Range("I62").Value = 2333.3333
$ws.Range("K62").Value = 2333.3333
$ws.Range("M62").Value = -1709.3333

$ws.Range("H65").Value = 3499.6667
$ws.Range("I65").Value = 2333.3333
$ws.Range("K65").Value = 11666.6665
$ws.Range("M65").Value = -8546.666499999999

$ws.Range("H96").Value = 1728.7
$ws.Range("J96").Value = 2800
$ws.Range("L96").Value = 2800
$ws.Range("N96").Value = -5546

$ws.Range("H100").Value = 6971716
$ws.Range("I100").Value = 11617228
$ws.Range("K100").Value = 23234456
$ws.Range("M100").Value = -23233915
